$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reset all existing offset values (column C, rows 1-30) to 0
$ws.Range("C1:C30").Value = 0

# 2. Add the new stations (rows 31-33)
$ws.Cells.Item(31,1).Value = "USACE"
$ws.Cells.Item(31,2).NumberFormat = "@"
$ws.Cells.Item(31,2).Value = "01480"
$ws.Cells.Item(31,3).NumberFormat = "0.00"
$ws.Cells.Item(31,3).Value = 0

$ws.Cells.Item(32,1).Value = "USACE"
$ws.Cells.Item(32,2).Value = 76560
$ws.Cells.Item(32,2).NumberFormat = "@"
$ws.Cells.Item(32,3).NumberFormat = "0.00"
$ws.Cells.Item(32,3).Value = 0

$ws.Cells.Item(33,1).Value = "USGS"
$ws.Cells.Item(33,2).NumberFormat = "@"
$ws.Cells.Item(33,2).Value = "073814675"
$ws.Cells.Item(33,3).NumberFormat = "0.00"
$ws.Cells.Item(33,3).Value = 0

# 3. Rename the station at row 7 (do this last so a fresh shared string is appended
#    rather than mutating the slot now reused by row 31's "01480")
$ws.Cells.Item(7,2).Value = "82742"

# 4. Update the sheet view to match the new extent / selection
[void]$ws.Range("C1:C33").Select()
